$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for rows 2-5 (new sensor readings) ---
$ws.Range("A2").Value = 45081.50694444445
$ws.Range("B2").Value = 23.541
$ws.Range("C2").Value = 16.304
$ws.Range("D2").Value = 4.244
$ws.Range("E2").Value = 49.627
$ws.Range("F2").Value = 41.058
$ws.Range("G2").Value = 18.526
$ws.Range("H2").Value = 61.718
$ws.Range("I2").Value = 28.505
$ws.Range("J2").Value = 12.158
$ws.Range("K2").Value = 18.744
$ws.Range("L2").Value = 19.59
$ws.Range("M2").Value = 20.45
$ws.Range("N2").Value = 5.915
$ws.Range("O2").Value = 18.422
$ws.Range("P2").Value = 25.954
$ws.Range("Q2").Value = 15.331
$ws.Range("R2").Value = 3.814
$ws.Range("S2").Value = 2.604
$ws.Range("T2").Value = 273.066
$ws.Range("U2").Value = 51.313
$ws.Range("V2").Value = 17.004
$ws.Range("W2").Value = 34.141
$ws.Range("X2").Value = 17.76
$ws.Range("Y2").Value = 2.295
$ws.Range("Z2").Value = 30.618
$ws.Range("AA2").Value = 15.02
$ws.Range("AB2").Value = 13.452
$ws.Range("AC2").Value = 15.735
$ws.Range("AD2").Value = 20.304
$ws.Range("AE2").Value = 3.64
$ws.Range("AF2").Value = 54.593
$ws.Range("AG2").Value = 9.503
$ws.Range("AH2").Value = 21.259

$ws.Range("A3").Value = 45081.51388888889
$ws.Range("B3").Value = 14.413
$ws.Range("C3").Value = 10.15
$ws.Range("D3").Value = 1.749
$ws.Range("E3").Value = 30.73
$ws.Range("F3").Value = 25.451
$ws.Range("G3").Value = 11.343
$ws.Range("H3").Value = 45.921
$ws.Range("I3").Value = 17.452
$ws.Range("J3").Value = 7.532
$ws.Range("K3").Value = 11.443
$ws.Range("L3").Value = 12.294
$ws.Range("M3").Value = 12.831
$ws.Range("N3").Value = 3.624
$ws.Range("O3").Value = 11.279
$ws.Range("P3").Value = 15.92
$ws.Range("Q3").Value = 9.628
$ws.Range("R3").Value = 1.616
$ws.Range("S3").Value = 1.011
$ws.Range("T3").Value = 164.365
$ws.Range("U3").Value = 31.621
$ws.Range("V3").Value = 10.411
$ws.Range("W3").Value = 20.993
$ws.Range("X3").Value = 11.199
$ws.Range("Y3").Value = 1.377
$ws.Range("Z3").Value = 21.776
$ws.Range("AA3").Value = 9.196
$ws.Range("AB3").Value = 8.321999999999999
$ws.Range("AC3").Value = 9.739000000000001
$ws.Range("AD3").Value = 12.819
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 41.526
$ws.Range("AG3").Value = 5.787
$ws.Range("AH3").Value = 13.016

$ws.Range("A4").Value = 45081.52083333334
$ws.Range("B4").Value = 13.452
$ws.Range("C4").Value = 9.657999999999999
$ws.Range("D4").Value = 1.224
$ws.Range("E4").Value = 28.857
$ws.Range("F4").Value = 23.897
$ws.Range("G4").Value = 10.587
$ws.Range("H4").Value = 42.061
$ws.Range("I4").Value = 16.288
$ws.Range("J4").Value = 7.109
$ws.Range("K4").Value = 10.721
$ws.Range("L4").Value = 11.608
$ws.Range("M4").Value = 12.118
$ws.Range("N4").Value = 3.382
$ws.Range("O4").Value = 10.527
$ws.Range("P4").Value = 14.894
$ws.Range("Q4").Value = 8.978999999999999
$ws.Range("R4").Value = 1.063
$ws.Range("S4").Value = 0.726
$ws.Range("T4").Value = 152.909
$ws.Range("U4").Value = 29.462
$ws.Range("V4").Value = 9.717000000000001
$ws.Range("W4").Value = 19.636
$ws.Range("X4").Value = 10.516
$ws.Range("Y4").Value = 1.292
$ws.Range("Z4").Value = 20.021
$ws.Range("AA4").Value = 8.583
$ws.Range("AB4").Value = 7.718
$ws.Range("AC4").Value = 9.045999999999999
$ws.Range("AD4").Value = 12.138
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 37.932
$ws.Range("AG4").Value = 5.419
$ws.Range("AH4").Value = 12.148

$ws.Range("A5").Value = 45081.52777777778
$ws.Range("B5").Value = 7.21
$ws.Range("C5").Value = 5.08
$ws.Range("D5").Value = 0.8
$ws.Range("E5").Value = 15.37
$ws.Range("F5").Value = 12.74
$ws.Range("G5").Value = 5.67
$ws.Range("H5").Value = 24.92
$ws.Range("I5").Value = 8.73
$ws.Range("J5").Value = 3.78
$ws.Range("K5").Value = 5.68
$ws.Range("L5").Value = 6.22
$ws.Range("M5").Value = 6.43
$ws.Range("N5").Value = 1.81
$ws.Range("O5").Value = 5.64
$ws.Range("P5").Value = 7.96
$ws.Range("Q5").Value = 4.89
$ws.Range("R5").Value = 0.76
$ws.Range("S5").Value = 0.44
$ws.Range("T5").Value = 78.53
$ws.Range("U5").Value = 15.89
$ws.Range("V5").Value = 5.21
$ws.Range("W5").Value = 10.52
$ws.Range("X5").Value = 5.68
$ws.Range("Y5").Value = 0.66
$ws.Range("Z5").Value = 11.58
$ws.Range("AA5").Value = 4.6
$ws.Range("AB5").Value = 4.18
$ws.Range("AC5").Value = 4.89
$ws.Range("AD5").Value = 6.49
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 22.6
$ws.Range("AG5").Value = 2.87
$ws.Range("AH5").Value = 6.51

# --- Delete row 6 entirely (dimension shrinks from AH6 to AH5) ---
$ws.Rows(6).Delete()

# --- Update column widths (ColumnWidth = stored_width - 5/6) ---
$ws.Columns(2).ColumnWidth = 7.166666666666667
$ws.Columns(3).ColumnWidth = 7.166666666666667
$ws.Columns(5).ColumnWidth = 7.166666666666667
$ws.Columns(7).ColumnWidth = 7.166666666666667
$ws.Columns(9).ColumnWidth = 7.166666666666667
$ws.Columns(10).ColumnWidth = 7.166666666666667
$ws.Columns(11).ColumnWidth = 7.166666666666667
$ws.Columns(12).ColumnWidth = 7.166666666666667
$ws.Columns(13).ColumnWidth = 7.166666666666667
$ws.Columns(15).ColumnWidth = 7.166666666666667
$ws.Columns(16).ColumnWidth = 7.166666666666667
$ws.Columns(17).ColumnWidth = 7.166666666666667
$ws.Columns(20).ColumnWidth = 8.166666666666666
$ws.Columns(22).ColumnWidth = 7.166666666666667
$ws.Columns(24).ColumnWidth = 7.166666666666667
$ws.Columns(26).ColumnWidth = 7.166666666666667
$ws.Columns(28).ColumnWidth = 7.166666666666667
$ws.Columns(29).ColumnWidth = 7.166666666666667
$ws.Columns(30).ColumnWidth = 7.166666666666667
$ws.Columns(34).ColumnWidth = 7.166666666666667
